$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: new header row for the "Number of employees / Assets / Turnover" table,
# styled the same as the existing "title" header rows (e.g. row 11 / row 19).
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B23:D23").Style = "title"

# Row 24: Micro
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

# Row 25: Small
$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""

# Row 26: Medium
$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "<250"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""

# Row 27: Large
$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">249"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
